$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.001.91'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.592.75'
$ws.Range('E3').Value = '  +7.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.599'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.68%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.577'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +14.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.03'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.21'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0840'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.13'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +15.30%  '
$ws.Range('D14').Value = '2.991.40'
$ws.Range('E14').Value = '  +7.35%  '
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('D16').Value = '2.592.82'
$ws.Range('E16').Value = '  +7.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.921'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +9.66%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.87%  '
$ws.Range('D19').Value = '46.167.40'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000101'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.15'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '271.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '29.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +39.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.56'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E31').Value = '  +4.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '38.80'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +14.98%  '
$ws.Range('E34').Value = '  -5.53%  '
$ws.Range('E35').Value = '  +2.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0837'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.19'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '149.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.121'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.122'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +42.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.89'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0328'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.07'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.54%  '
$ws.Range('D46').Value = '2.152.54'
$ws.Range('E46').Value = '  +8.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.997'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '93.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '108.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.50%  '
